$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 12.79549166666667
$ws.Cells.Item(2, 8).Value = 38.386475
$ws.Cells.Item(2, 9).Value = 0.5145949251267348
$ws.Cells.Item(2, 10).Value = 0.5145949251267348
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.1743356666666667
$ws.Cells.Item(2, 14).Value = 0.523007
$ws.Cells.Item(2, 15).Value = 0.02303006925000699
$ws.Cells.Item(2, 16).Value = 0.02303006925000699
$ws.Cells.Item(2, 17).Value = 2.230710570036111
$ws.Cells.Item(2, 18).Value = 20.076395130325
$ws.Cells.Item(2, 19).Value = 0.01185115676137086
$ws.Cells.Item(2, 20).Value = 0.01185115676137086

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 12.79549166666667
$ws.Cells.Item(3, 8).Value = 38.386475
$ws.Cells.Item(3, 9).Value = 0.5145949251267348
$ws.Cells.Item(3, 10).Value = 0.5145949251267348
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 2.401382333333334
$ws.Cells.Item(3, 14).Value = 7.204147000000001
$ws.Cells.Item(3, 15).Value = 0.3172271198994089
$ws.Cells.Item(3, 16).Value = 0.3172271198994089
$ws.Cells.Item(3, 17).Value = 30.72686763464722
$ws.Cells.Item(3, 18).Value = 276.541808711825
$ws.Cells.Item(3, 19).Value = 0.163243466012806
$ws.Cells.Item(3, 20).Value = 0.163243466012806

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 12.79549166666667
$ws.Cells.Item(4, 8).Value = 38.386475
$ws.Cells.Item(4, 9).Value = 0.5145949251267348
$ws.Cells.Item(4, 10).Value = 0.5145949251267348
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 4.994197
$ws.Cells.Item(4, 14).Value = 14.982591
$ws.Cells.Item(4, 15).Value = 0.6597428108505842
$ws.Cells.Item(4, 16).Value = 0.6597428108505842
$ws.Cells.Item(4, 17).Value = 63.90320609519166
$ws.Cells.Item(4, 18).Value = 575.128854856725
$ws.Cells.Item(4, 19).Value = 0.3395003023525579
$ws.Cells.Item(4, 20).Value = 0.3395003023525579

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 9.803896
$ws.Cells.Item(5, 8).Value = 29.411688
$ws.Cells.Item(5, 9).Value = 0.3942822409249843
$ws.Cells.Item(5, 10).Value = 0.3942822409249843
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.1743356666666667
$ws.Cells.Item(5, 14).Value = 0.523007
$ws.Cells.Item(5, 15).Value = 0.02303006925000699
$ws.Cells.Item(5, 16).Value = 0.02303006925000699
$ws.Cells.Item(5, 17).Value = 1.709168745090667
$ws.Cells.Item(5, 18).Value = 15.382518705816
$ws.Cells.Item(5, 19).Value = 0.009080347312550328
$ws.Cells.Item(5, 20).Value = 0.009080347312550328

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 9.803896
$ws.Cells.Item(6, 8).Value = 29.411688
$ws.Cells.Item(6, 9).Value = 0.3942822409249843
$ws.Cells.Item(6, 10).Value = 0.3942822409249843
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 2.401382333333334
$ws.Cells.Item(6, 14).Value = 7.204147000000001
$ws.Cells.Item(6, 15).Value = 0.3172271198994089
$ws.Cells.Item(6, 16).Value = 0.3172271198994089
$ws.Cells.Item(6, 17).Value = 23.54290265223734
$ws.Cells.Item(6, 18).Value = 211.886123870136
$ws.Cells.Item(6, 19).Value = 0.1250770197161176
$ws.Cells.Item(6, 20).Value = 0.1250770197161176

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 9.803896
$ws.Cells.Item(7, 8).Value = 29.411688
$ws.Cells.Item(7, 9).Value = 0.3942822409249843
$ws.Cells.Item(7, 10).Value = 0.3942822409249843
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 4.994197
$ws.Cells.Item(7, 14).Value = 14.982591
$ws.Cells.Item(7, 15).Value = 0.6597428108505842
$ws.Cells.Item(7, 16).Value = 0.6597428108505842
$ws.Cells.Item(7, 17).Value = 48.962587991512
$ws.Cells.Item(7, 18).Value = 440.663291923608
$ws.Cells.Item(7, 19).Value = 0.2601248738963163
$ws.Cells.Item(7, 20).Value = 0.2601248738963163

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 2.265785
$ws.Cells.Item(8, 8).Value = 6.797355
$ws.Cells.Item(8, 9).Value = 0.09112283394828093
$ws.Cells.Item(8, 10).Value = 0.09112283394828093
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 0.1743356666666667
$ws.Cells.Item(8, 14).Value = 0.523007
$ws.Cells.Item(8, 15).Value = 0.02303006925000699
$ws.Cells.Item(8, 16).Value = 0.02303006925000699
$ws.Cells.Item(8, 17).Value = 0.3950071384983334
$ws.Cells.Item(8, 18).Value = 3.555064246485
$ws.Cells.Item(8, 19).Value = 0.002098565176085798
$ws.Cells.Item(8, 20).Value = 0.002098565176085798

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 2.265785
$ws.Cells.Item(9, 8).Value = 6.797355
$ws.Cells.Item(9, 9).Value = 0.09112283394828093
$ws.Cells.Item(9, 10).Value = 0.09112283394828093
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 2.401382333333334
$ws.Cells.Item(9, 14).Value = 7.204147000000001
$ws.Cells.Item(9, 15).Value = 0.3172271198994089
$ws.Cells.Item(9, 16).Value = 0.3172271198994089
$ws.Cells.Item(9, 17).Value = 5.441016070131668
$ws.Cells.Item(9, 18).Value = 48.96914463118501
$ws.Cells.Item(9, 19).Value = 0.02890663417048524
$ws.Cells.Item(9, 20).Value = 0.02890663417048524

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 2.265785
$ws.Cells.Item(10, 8).Value = 6.797355
$ws.Cells.Item(10, 9).Value = 0.09112283394828093
$ws.Cells.Item(10, 10).Value = 0.09112283394828093
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 4.994197
$ws.Cells.Item(10, 14).Value = 14.982591
$ws.Cells.Item(10, 15).Value = 0.6597428108505842
$ws.Cells.Item(10, 16).Value = 0.6597428108505842
$ws.Cells.Item(10, 17).Value = 11.315776649645
$ws.Cells.Item(10, 18).Value = 101.841989846805
$ws.Cells.Item(10, 19).Value = 0.06011763460170989
$ws.Cells.Item(10, 20).Value = 0.06011763460170989

